# Background section finished - add the new work-log entry to the log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 30: activity description + hours spent (2 hrs)
$ws.Range("B30").Value = "Finished updating Backgrounds section"
$ws.Range("C30").Value = 2

# Scroll the window down a bit and land the selection on D30, like a user who
# just finished typing the new row would leave the sheet.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D30").Select()
